$d = $word.ActiveDocument
$t = $d.Tables(1)
$nl = [char]11

$t.Cell(1, 1).Range.Text = "92 x 70" + $nl + "  7    0" + $nl + "  ----" + $nl + "9|    |" + $nl + "2|    |"
$t.Cell(1, 2).Range.Text = "93 x 84" + $nl + "  8    4" + $nl + "  ----" + $nl + "9|    |" + $nl + "3|    |"
$t.Cell(1, 3).Range.Text = "30 x 12" + $nl + "  1    2" + $nl + "  ----" + $nl + "3|    |" + $nl + "0|    |"
$t.Cell(2, 1).Range.Text = "13 x 31" + $nl + "  3    1" + $nl + "  ----" + $nl + "1|    |" + $nl + "3|    |"
$t.Cell(2, 2).Range.Text = "58 x 31" + $nl + "  3    1" + $nl + "  ----" + $nl + "5|    |" + $nl + "8|    |"
$t.Cell(2, 3).Range.Text = "12 x 15" + $nl + "  1    5" + $nl + "  ----" + $nl + "1|    |" + $nl + "2|    |"
$t.Cell(3, 1).Range.Text = "53 x 96" + $nl + "  9    6" + $nl + "  ----" + $nl + "5|    |" + $nl + "3|    |"
$t.Cell(3, 2).Range.Text = "27 x 64" + $nl + "  6    4" + $nl + "  ----" + $nl + "2|    |" + $nl + "7|    |"
$t.Cell(3, 3).Range.Text = "55 x 68" + $nl + "  6    8" + $nl + "  ----" + $nl + "5|    |" + $nl + "5|    |"
$t.Cell(4, 1).Range.Text = "78 x 51" + $nl + "  5    1" + $nl + "  ----" + $nl + "7|    |" + $nl + "8|    |"
$t.Cell(4, 2).Range.Text = "74 x 47" + $nl + "  4    7" + $nl + "  ----" + $nl + "7|    |" + $nl + "4|    |"
$t.Cell(4, 3).Range.Text = "65 x 44" + $nl + "  4    4" + $nl + "  ----" + $nl + "6|    |" + $nl + "5|    |"
$t.Cell(5, 1).Range.Text = "75 x 38" + $nl + "  3    8" + $nl + "  ----" + $nl + "7|    |" + $nl + "5|    |"
$t.Cell(5, 2).Range.Text = "98 x 72" + $nl + "  7    2" + $nl + "  ----" + $nl + "9|    |" + $nl + "8|    |"
$t.Cell(5, 3).Range.Text = "30 x 85" + $nl + "  8    5" + $nl + "  ----" + $nl + "3|    |" + $nl + "0|    |"
